$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "<field>_old" -> "<field>_FV2404", "<field>_new" -> "<field>_FV2410"
# Columns A..J and L..U hold the "_old"/"_new" suffixed headers; column K is "diff" (unchanged).
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string]) {
        if ($val.EndsWith("_old")) {
            $base = $val.Substring(0, $val.Length - 4)
            $cell.Value2 = $base + "_FV2404"
        } elseif ($val.EndsWith("_new")) {
            $base = $val.Substring(0, $val.Length - 4)
            $cell.Value2 = $base + "_FV2410"
        }
    }
}

# 2. Freeze the top row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Convert the data range into a table (ListObject)
$range = $ws.Range("A1:U80")
$list = $ws.ListObjects.Add(1, $range, $null, 1)
$list.Name = "Table1"
